$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45264
$ws.Range("M2").Value = 150
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = '$/caja 10 kilos'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 1500

# Row 3
$ws.Range("D3").Value = 45251
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 150
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("Q3").Value = '$/bandeja 10 kilos'
$ws.Range("R3").Value = 'Provincia de Limarí'
$ws.Range("S3").Value = 2000

# Row 4
$ws.Range("D4").Value = 45257
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("Q4").Value = '$/bandeja 10 kilos'
$ws.Range("R4").Value = 'Provincia de Limarí'
$ws.Range("S4").Value = 2000

# Row 5
$ws.Range("D5").Value = 44505
$ws.Range("K5").Value = 'Californiana(o)'
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 1500

# Row 6
$ws.Range("K6").Value = 'Golden Nugget'
$ws.Range("M6").Value = 50

# Row 7
$ws.Range("D7").Value = 44902
$ws.Range("L7").Value = 'Especial'
$ws.Range("M7").Value = 60
$ws.Range("Q7").Value = '$/caja 10 kilos'

# Row 8
$ws.Range("D8").Value = 44902
$ws.Range("M8").Value = 70
$ws.Range("N8").Value = 13000
$ws.Range("O8").Value = 13000
$ws.Range("P8").Value = 13000
$ws.Range("S8").Value = 1300
